$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing the existing rows 2 and 3 down to 3 and 4.
$ws.Rows.Item(2).Insert()

# Excel copies the bold/centered header formatting onto the newly inserted row;
# drop it so the row matches the other plain data rows.
$ws.Rows.Item(2).ClearFormats()

# The date column on the new row needs the same date/time number format as the rest
# of the Fecha column.
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the updated weekly record.
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 45020
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100104
$ws.Cells.Item(2, 8).Value = "Frutos de pepita"
$ws.Cells.Item(2, 9).Value = 100104003
$ws.Cells.Item(2, 10).Value = "Membrillo"
$ws.Cells.Item(2, 11).Value = "Champion"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 12000
$ws.Cells.Item(2, 15).Value = 12000
$ws.Cells.Item(2, 16).Value = 12000
$ws.Cells.Item(2, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(2, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(2, 19).Value = 667
$ws.Cells.Item(2, 20).Value = 18
